# Incorporacion de Logica para Join con Efectores solapa SIF-SIGEHOS
#
# Adds a new "EfectorSigehos" column (D) to the Tabla2 table on Hoja1,
# populated with the accent-stripped / SIGEHOS-normalized version of the
# existing "EfectorObjetivos" column (C), and updates the sheet selection
# and column width accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the new table column (grows the table range from A1:C38 to A1:D38
# and the AutoFilter along with it).
$newCol = $lo.ListColumns.Add()

# Header
$ws.Range("D1").Value = "EfectorSigehos"

# Data rows: EfectorSigehos values (SIGEHOS-side effector names), one per
# table row (rows 2-38 on the sheet).
$values = @(
    "Pirovano",
    "Santojanni",
    "Fernandez",
    "Alvarez",
    "Quemados",
    "Rocca",
    "Santa Lucia",
    "Ferrer",
    "Penna",
    "Muñiz",
    "Lagleyze",
    "Argerich",
    "Rivadavia",
    "Udaondo",
    "Velez Sarsfield",
    "Durand",
    "Gutierrez",
    "Tornu",
    "Curie",
    "Zubizarreta",
    "Piñero",
    "Sarda",
    "Elizalde",
    "Ramos Mejia",
    "Borda",
    "Moyano",
    "IREP",
    "Dueñas",
    "Tobar Garcia",
    "Quinquela Martin",
    "Talleres Protegidos",
    "Alvear",
    "Carrillo",
    "SAME",
    "Grierson",
    "Barrio 31",
    "Turismo"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Widen the new column like the other data columns.
$ws.Columns.Item(4).ColumnWidth = 18.7109375

# Move the active selection, as captured in the saved workbook.
$ws.Range("B8").Select()
